$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Meeting with Hrn. Rudin: the "R5 / Apache-Tomcat" risk row is removed
# from the risk table entirely; the rows below it shift up by one.
$ws.Rows("12").Delete()

# Restore the shared formula grouping for the "max.Schaden x Wahrscheinlichkeit"
# column (F8:F12) that the row deletion left as individual formulas.
$ws.Range("F8:F12").Formula = "=D8*E8"

# Update the window's scroll position / active selection to reflect the
# state captured after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C16").Select()
